$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reduce the tab-ratio (width of sheet-tabs area vs horizontal scrollbar) from 60% to 50%
$wb.Windows.Item(1).TabRatio = 0.5

# Fill in the previously-blank "payment_term_id" column (H) for the new test rows.
# H4 / H8 reuse the same look as the existing H1/H2 cells (Test1's payment term).
$h4 = $ws.Range("H4")
$h4.Value = "z0bug.payment_1"
$h4.Font.Size = 9
$h4.Font.Name = "arial"
$h4.Font.Color = 0
$h4.Font.Family = 2

$h5 = $ws.Range("H5")
$h5.Value = "z0bug.payment_5"
$h5.Font.Size = 9
$h5.Font.Name = "arial"
$h5.Font.Color = 0
$h5.Font.Family = 0

$h7 = $ws.Range("H7")
$h7.Value = "z0bug.payment_4"
$h7.Font.Size = 9
$h7.Font.Name = "arial"
$h7.Font.Color = 0
$h7.Font.Family = 0

$h8 = $ws.Range("H8")
$h8.Value = "z0bug.payment_1"
$h8.Font.Size = 9
$h8.Font.Name = "arial"
$h8.Font.Color = 0
$h8.Font.Family = 2

# Move / leave the active selection on H7, matching the saved cursor position.
$ws.Range("H7").Select()
